# Test Data Added for Slovakia market
# - Append new device codes to the shared "NGC-3475/..." value (affects the
#   B4 cell on every sheet, since they all reference the same shared string).
# - Update the remembered selection on a few sheets.
# - Leave "Other Zetfas" as the active sheet/tab when the workbook is saved.

$wb = $excel.ActiveWorkbook

# Update the shared device-code text. All seven sheets' B4 cell points at
# the same shared string, so touch every sheet to keep them sharing one
# (updated) string table entry instead of forking a new one.
$deviceCodes = "NGC-3475/T1832/1858/NGC-2930/T3173/T1832"
foreach ($sheetName in @("Detectors_STI", "Call points STI", "Ancillary STI", "Other STI", "Detectors_Zetfas", "Ancillary Zetfas", "Other Zetfas")) {
    $wb.Worksheets.Item($sheetName).Range("B4").Value = $deviceCodes
}

# Ancillary STI: selection moves from B29 to B4.
$wsAncillarySTI = $wb.Worksheets.Item("Ancillary STI")
$wsAncillarySTI.Activate() | Out-Null
$wsAncillarySTI.Range("B4").Select() | Out-Null

# Other STI: selection collapses from the B4 / C6:C21 multi-range to B4.
$wsOtherSTI = $wb.Worksheets.Item("Other STI")
$wsOtherSTI.Activate() | Out-Null
$wsOtherSTI.Range("B4").Select() | Out-Null

# Detectors_STI: selection moves from A29 to C12, and it's no longer the
# active/visible tab once we finish on "Other Zetfas" below.
$wsDetectorsSTI = $wb.Worksheets.Item("Detectors_STI")
$wsDetectorsSTI.Activate() | Out-Null
$wsDetectorsSTI.Range("C12").Select() | Out-Null

# Other Zetfas: becomes the active tab, selection moves from B7 to B4. This
# is the last sheet touched so it ends up the active tab on save.
$wsOtherZetfas = $wb.Worksheets.Item("Other Zetfas")
$wsOtherZetfas.Activate() | Out-Null
$wsOtherZetfas.Range("B4").Select() | Out-Null
